$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary updates ---
$ws.Range("E11").Value = 924000
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 2

# --- Insert a new row above the last worker row (old row 19) so the
#     table grows from 4 workers to 5. Copy formatting from the row
#     above (row 18) so the new row matches the other data rows. ---
$ws.Rows("19").Insert()
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Re-write the worker table (rows 16-19 are now the "normal" style
#     rows, row 20 keeps the old bottom-border style row that got
#     pushed down by the insert). Order changed: GINA now comes first. ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "43924603"
$ws.Range("D16").Value = "GINA PAOLA CARRASQUILLA ESQUIVEL"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 100000
$ws.Range("G16").Value = 2500000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143373133"
$ws.Range("D17").Value = "ALEXANDER PATRON GONZALEZ"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 64000
$ws.Range("G17").Value = 1600000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047484625"
$ws.Range("D18").Value = "CRISTIAN MOGOLLON HOYOS"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 64000
$ws.Range("G18").Value = 1600000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1098753514"
$ws.Range("D19").Value = "MARIA CAMILA SILVA DIAZ"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 64000
$ws.Range("G19").Value = 1600000

# --- New row 20 (previously row 19, now shifted down) becomes a new
#     NIT-level summary line for the account itself. ---
$ws.Range("B20").Value = "NIT"
$ws.Range("C20").Value = "9009513862"
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = "1606"
$ws.Range("F20").Value = 632000
$ws.Range("G20").Value = 0
